$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Seccion" column for the remaining rows of each group
$ws.Range("A3").Value = "A-1"
$ws.Range("A4").Value = "A-1"
$ws.Range("A6").Value = "A-2"
$ws.Range("A7").Value = "A-2"

# Update the active cell selection
$ws.Range("E8").Select()
